$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SO_App now appends the newly submitted Special Order right after
# creation, leaving a blank spacer row between the previous last entry
# (row 2) and the freshly submitted one (row 4).

# --- Row 2 cleanup: the trailing blank "ADDRESS" placeholder cell is gone ---
$ws.Cells.Item(2, 13).ClearContents()

# --- Row 3: blank spacer row (empty text in every column, A:P) ---
$blankCols = 1..16
foreach ($col in $blankCols) {
    $ws.Cells.Item(3, $col).Value = "'"
}

# --- Row 4: newly submitted Special Order ---
$ws.Cells.Item(4, 1).Value = "'02/19/2024"
$ws.Cells.Item(4, 2).Value = "SO240219001"
$ws.Cells.Item(4, 3).Value = "ab"
$ws.Cells.Item(4, 4).Value = "'7894561230"
$ws.Cells.Item(4, 5).Value = "a"
$ws.Cells.Item(4, 6).Value = "a"
$ws.Cells.Item(4, 7).Value = 12
$ws.Cells.Item(4, 8).Value = 34
$ws.Cells.Item(4, 9).Value = "AMA"
$ws.Cells.Item(4, 10).Value = "OTHER"
$ws.Cells.Item(4, 11).Value = "ab"
$ws.Cells.Item(4, 12).Value = "NO"
$ws.Cells.Item(4, 13).Value = "'"
